$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7647.4
$ws.Range("J43").Value = 1004.5714
$ws.Range("L43").Value = 1004.5714
$ws.Range("N43").Value = -1142.5714

$ws.Range("H127").Value = 2089.6086
$ws.Range("I127").Value = 1025
$ws.Range("J127").Value = 2657.4
$ws.Range("K127").Value = 3075
$ws.Range("L127").Value = 7972.200000000001
$ws.Range("M127").Value = 1885
$ws.Range("N127").Value = -17892.2

$ws.Range("H141").Value = 4384.032
$ws.Range("I141").Value = 4392.857
$ws.Range("K141").Value = 13178.571
$ws.Range("M141").Value = -7998.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2608.6538
$ws.Range("I32").Value = 2307.8767
$ws.Range("J32").Value = 7000
$ws.Range("K32").Value = 2307.8767
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = -2020.8767
$ws.Range("N32").Value = -7574

$ws.Range("H45").Value = 1140.6428
$ws.Range("I45").Value = 629.75
$ws.Range("J45").Value = 1345
$ws.Range("K45").Value = 629.75
$ws.Range("L45").Value = 1345
$ws.Range("M45").Value = -252.75
$ws.Range("N45").Value = -2099

$ws.Range("H63").Value = 3033282.5
$ws.Range("I63").Value = 8334433
$ws.Range("J63").Value = 4053.5715
$ws.Range("K63").Value = 8334433
$ws.Range("L63").Value = 4053.5715
$ws.Range("M63").Value = -8333747
$ws.Range("N63").Value = -5425.5715

$ws.Range("H66").Value = 3033282.5
$ws.Range("I66").Value = 8334433
$ws.Range("J66").Value = 4053.5715
$ws.Range("K66").Value = 41672165
$ws.Range("L66").Value = 20267.8575
$ws.Range("M66").Value = -41668733
$ws.Range("N66").Value = -27131.8575

$ws.Range("H102").Value = 100002140
$ws.Range("I102").Value = 125002200
$ws.Range("K102").Value = 125002200
$ws.Range("M102").Value = -125000578

$ws.Range("H110").Value = 815.0625
$ws.Range("I110").Value = 728.8570999999999
$ws.Range("J110").Value = 882.1111
$ws.Range("K110").Value = 728.8570999999999
$ws.Range("L110").Value = 882.1111
$ws.Range("M110").Value = 1316.1429
$ws.Range("N110").Value = -4972.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 760
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = -513
$ws.Range("N16").Value = -1174

$ws.Range("H99").Value = 62502110
$ws.Range("I99").Value = 76925310
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 76925310
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = -76923812
$ws.Range("N99").Value = -4596

$ws.Range("H113").Value = 760
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -4940

$ws.Range("H126").Value = 62502110
$ws.Range("I126").Value = 76925310
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 230775930
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -230773460
$ws.Range("N126").Value = -9740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 115
$ws.Range("J10").Value = 150
$ws.Range("L10").Value = 450
$ws.Range("N10").Value = -728

$ws.Range("H44").Value = 949.2857
$ws.Range("I44").Value = 649.625
$ws.Range("J44").Value = 1348.8334
$ws.Range("K44").Value = 1948.875
$ws.Range("L44").Value = 4046.5002
$ws.Range("M44").Value = -1550.875
$ws.Range("N44").Value = -4842.5002

$ws.Range("H80").Value = 3868.8096
$ws.Range("I80").Value = 1563.75
$ws.Range("J80").Value = 4411.1763
$ws.Range("K80").Value = 4691.25
$ws.Range("L80").Value = 13233.5289
$ws.Range("M80").Value = -3755.25
$ws.Range("N80").Value = -15105.5289

$ws.Range("H83").Value = 3868.8096
$ws.Range("I83").Value = 1563.75
$ws.Range("J83").Value = 4411.1763
$ws.Range("K83").Value = 14073.75
$ws.Range("L83").Value = 39700.5867
$ws.Range("M83").Value = -9393.75
$ws.Range("N83").Value = -49060.5867

$ws.Range("H92").Value = 899.5
$ws.Range("I92").Value = 818.8570999999999
$ws.Range("J92").Value = 962.2222
$ws.Range("K92").Value = 2456.5713
$ws.Range("L92").Value = 2886.6666
$ws.Range("M92").Value = -1208.5713
$ws.Range("N92").Value = -5382.6666

$ws.Range("H104").Value = 6465
$ws.Range("I104").Value = 2097.5
$ws.Range("J104").Value = 7629.6665
$ws.Range("K104").Value = 6292.5
$ws.Range("L104").Value = 22888.9995
$ws.Range("M104").Value = -3671.5
$ws.Range("N104").Value = -28130.9995

$ws.Range("H134").Value = 3254.6667
$ws.Range("I134").Value = 3377.5
$ws.Range("J134").Value = 3114.2856
$ws.Range("K134").Value = 10132.5
$ws.Range("L134").Value = 9342.856800000001
$ws.Range("M134").Value = -5062.5
$ws.Range("N134").Value = -19482.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1537.7587
$ws.Range("I7").Value = 1305.1052
$ws.Range("J7").Value = 1979.8
$ws.Range("K7").Value = 1305.1052
$ws.Range("L7").Value = 1979.8
$ws.Range("M7").Value = -1193.1052
$ws.Range("N7").Value = -2203.8

$ws.Range("H16").Value = 14410.7
$ws.Range("I16").Value = 774.75
$ws.Range("J16").Value = 23501.334
$ws.Range("K16").Value = 774.75
$ws.Range("L16").Value = 23501.334
$ws.Range("M16").Value = -604.75
$ws.Range("N16").Value = -23841.334

$ws.Range("H40").Value = 1900
$ws.Range("I40").Value = 1925
$ws.Range("J40").Value = 1700
$ws.Range("K40").Value = 1925
$ws.Range("L40").Value = 1700
$ws.Range("M40").Value = -1789
$ws.Range("N40").Value = -1972

$ws.Range("H122").Value = 125002400
$ws.Range("I122").Value = 142858720
$ws.Range("J122").Value = 83337660
$ws.Range("K122").Value = 428576160
$ws.Range("L122").Value = 250012980
$ws.Range("M122").Value = -428573710
$ws.Range("N122").Value = -250017880

$ws.Range("H126").Value = 1537.7587
$ws.Range("I126").Value = 1305.1052
$ws.Range("J126").Value = 1979.8
$ws.Range("K126").Value = 3915.3156
$ws.Range("L126").Value = 5939.4
$ws.Range("M126").Value = -1445.3156
$ws.Range("N126").Value = -10879.4
